$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the existing "Monthly" sheet to "MonthlyLots" (content/sheetId
#    stay put - only the name + a couple of view properties change).
# ---------------------------------------------------------------------------
$wsLots = $wb.Worksheets.Item("Monthly")
$wsLots.Name = "MonthlyLots"

# view tweak on MonthlyLots: selection changes from a single active cell to a
# full-sheet-ish range selection (no explicit active cell).
$wsLots.Range("A1:XFD2").Select()

# ---------------------------------------------------------------------------
# 2. Insert two new sheets right after MonthlyLots: "Monthly" and "Calc"
# ---------------------------------------------------------------------------
$wsMonthly = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsLots)
$wsMonthly.Name = "Monthly"

$wsCalc = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsMonthly)
$wsCalc.Name = "Calc"

# ---------------------------------------------------------------------------
# 3. Populate the new "Monthly" sheet - header row + one royalty calc record
#    (mirrors the first data row of MonthlyLots).
# ---------------------------------------------------------------------------
$wsMonthly.Range("A1").Value = "Row"
$wsMonthly.Range("B1").Value = "ExtractMonth"
$wsMonthly.Range("C1").Value = "ProdMonth"
$wsMonthly.Range("D1").Value = "WellId"
$wsMonthly.Range("E1").Value = "Product"
$wsMonthly.Range("F1").Value = "AmendNo"
$wsMonthly.Range("G1").Value = "ProdHours"
$wsMonthly.Range("H1").Value = "ProdVol"
$wsMonthly.Range("I1").Value = "TransPrice"
$wsMonthly.Range("J1").Value = "WellHeadPrice"
$wsMonthly.Range("K1").Value = "TransRate"
$wsMonthly.Range("L1").Value = "ProcessingRate"

$wsMonthly.Range("A2").Value = 1
$wsMonthly.Range("B2").Value = 42276
$wsMonthly.Range("C2").Value = 201501
$wsMonthly.Range("D2").Value = 6
$wsMonthly.Range("E2").Value = "Oil"
$wsMonthly.Range("F2").Value = 2
$wsMonthly.Range("G2").Value = 740
$wsMonthly.Range("H2").Value = 100
$wsMonthly.Range("I2").Value = 2.2000000000000002
$wsMonthly.Range("J2").Value = 221.123456
$wsMonthly.Range("K2").Value = 2.1234549999999999
$wsMonthly.Range("L2").Value = 0.123455

# Re-use the existing date / quote-prefix formats already present on
# MonthlyLots (B2 = short date, D2 = quote-prefixed number) instead of
# minting brand-new number formats.
$wsLots.Range("B2").Copy()
$wsMonthly.Range("B2").PasteSpecial(-4122)
$wsLots.Range("D2").Copy()
$wsMonthly.Range("D2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsMonthly.Columns.Item(2).ColumnWidth = 12

$wsMonthly.Range("D27:I28").Select()
$wsMonthly.Application.ActiveWindow.RangeSelection.Item(1).Activate()

# ---------------------------------------------------------------------------
# 4. Populate the new "Calc" sheet - header row for a royalty calc record.
# ---------------------------------------------------------------------------
$wsCalc.Range("A1").Value = "ProdMonth"
$wsCalc.Range("B1").Value = "WellId"
$wsCalc.Range("C1").Value = "K"
$wsCalc.Range("D1").Value = "X"
$wsCalc.Range("E1").Value = "C"
$wsCalc.Range("F1").Value = "D"
$wsCalc.Range("G1").Value = "RoyaltyPrice"
$wsCalc.Range("H1").Value = "RoyaltyVolume"
$wsCalc.Range("I1").Value = "ProvCrownRoyaltyRate"
$wsCalc.Range("J1").Value = "ProvCrownUsedRoyaltyRate"
$wsCalc.Range("K1").Value = "IOGR1995RoyaltyRate"
$wsCalc.Range("L1").Value = "GorrRoyaltyRate"
$wsCalc.Range("M1").Value = "ProvCrownRoyaltyVolume"
$wsCalc.Range("N1").Value = "GorrRoyaltyVolume"
$wsCalc.Range("O1").Value = "IOGR1995RoyaltyVolume"
$wsCalc.Range("P1").Value = "ProvCrownRoyaltyValue"
$wsCalc.Range("Q1").Value = "IOGR1995RoyaltyValue"
$wsCalc.Range("R1").Value = "GorrRoyaltyValue"
$wsCalc.Range("S1").Value = "RoyaltyValuePreDeductions"
$wsCalc.Range("T1").Value = "RoyaltyTransportation"
$wsCalc.Range("U1").Value = "RoyaltyProcessing"
$wsCalc.Range("V1").Value = "RoyaltyDeductions"
$wsCalc.Range("W1").Value = "RoyaltyValue"
$wsCalc.Range("X1").Value = "CommencementPeriod"
$wsCalc.Range("Y1").Value = "Message"
$wsCalc.Range("Z1").Value = "GorrMessage"

$wsCalc.Application.ActiveWindow.ScrollColumn = 7
$wsCalc.Range("Z1").Select()
$wsCalc.Activate()

# ---------------------------------------------------------------------------
# 5. RoyaltyMaster: the tab-selected flag moves to Calc, and the selection
#    on RoyaltyMaster itself changes.
# ---------------------------------------------------------------------------
$wsRoyaltyMaster = $wb.Worksheets.Item("RoyaltyMaster")
$wsRoyaltyMaster.Range("D18").Select()

# ---------------------------------------------------------------------------
# 6. Book-level view: firstSheet/activeTab move forward now that two sheets
#    were inserted ahead of the previously-visible tab range.
# ---------------------------------------------------------------------------
$wb.Windows.Item(1).ScrollWorkbookTabs(2)

$wsCalc.Activate()
